$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.047456
$ws.Range("H2").Value = 0.142368
$ws.Range("I2").Value = 0.03340144944197188
$ws.Range("J2").Value = 0.03340144944197188
$ws.Range("M2").Value = 2.507757
$ws.Range("N2").Value = 7.523270999999999
$ws.Range("O2").Value = 0.07648041298707947
$ws.Range("P2").Value = 0.07648041298707947
$ws.Range("Q2").Value = 0.119008116192
$ws.Range("R2").Value = 1.071073045728
$ws.Range("S2").Value = 0.002554556647689065
$ws.Range("T2").Value = 0.002554556647689064

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.047456
$ws.Range("H3").Value = 0.142368
$ws.Range("I3").Value = 0.03340144944197188
$ws.Range("J3").Value = 0.03340144944197188
$ws.Range("O3").Value = 0.6219651214303167
$ws.Range("P3").Value = 0.6219651214303167
$ws.Range("Q3").Value = 0.9678150855573333
$ws.Range("R3").Value = 8.710335770016
$ws.Range("S3").Value = 0.02077453655812463
$ws.Range("T3").Value = 0.02077453655812462

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.047456
$ws.Range("H4").Value = 0.142368
$ws.Range("I4").Value = 0.03340144944197188
$ws.Range("J4").Value = 0.03340144944197188
$ws.Range("M4").Value = 9.887829999999999
$ws.Range("N4").Value = 29.66349
$ws.Range("O4").Value = 0.3015544655826039
$ws.Range("P4").Value = 0.301554465582604
$ws.Range("Q4").Value = 0.46923686048
$ws.Range("R4").Value = 4.22313174432
$ws.Range("S4").Value = 0.01007235623615819
$ws.Range("T4").Value = 0.01007235623615819

$ws.Range("I5").Value = 0.7427665150281654
$ws.Range("J5").Value = 0.7427665150281653
$ws.Range("M5").Value = 2.507757
$ws.Range("N5").Value = 7.523270999999999
$ws.Range("O5").Value = 0.07648041298707947
$ws.Range("P5").Value = 0.07648041298707947
$ws.Range("Q5").Value = 2.646449336804
$ws.Range("R5").Value = 23.818044031236
$ws.Range("S5").Value = 0.05680708982232786
$ws.Range("T5").Value = 0.05680708982232785

$ws.Range("I6").Value = 0.7427665150281654
$ws.Range("J6").Value = 0.7427665150281653
$ws.Range("O6").Value = 0.6219651214303167
$ws.Range("P6").Value = 0.6219651214303167
$ws.Range("S6").Value = 0.461974865713866
$ws.Range("T6").Value = 0.4619748657138659

$ws.Range("I7").Value = 0.7427665150281654
$ws.Range("J7").Value = 0.7427665150281653
$ws.Range("M7").Value = 9.887829999999999
$ws.Range("N7").Value = 29.66349
$ws.Range("O7").Value = 0.3015544655826039
$ws.Range("P7").Value = 0.301554465582604
$ws.Range("Q7").Value = 10.43467973409333
$ws.Range("R7").Value = 93.91211760684
$ws.Range("S7").Value = 0.2239845594919715
$ws.Range("T7").Value = 0.2239845594919715

$ws.Range("G8").Value = 0.3180153333333334
$ws.Range("H8").Value = 0.9540460000000001
$ws.Range("I8").Value = 0.2238320355298628
$ws.Range("J8").Value = 0.2238320355298628
$ws.Range("M8").Value = 2.507757
$ws.Range("N8").Value = 7.523270999999999
$ws.Range("O8").Value = 0.07648041298707947
$ws.Range("P8").Value = 0.07648041298707947
$ws.Range("Q8").Value = 0.797505178274
$ws.Range("R8").Value = 7.177546604465999
$ws.Range("S8").Value = 0.01711876651706255
$ws.Range("T8").Value = 0.01711876651706255

$ws.Range("G9").Value = 0.3180153333333334
$ws.Range("H9").Value = 0.9540460000000001
$ws.Range("I9").Value = 0.2238320355298628
$ws.Range("J9").Value = 0.2238320355298628
$ws.Range("O9").Value = 0.6219651214303167
$ws.Range("P9").Value = 0.6219651214303167
$ws.Range("Q9").Value = 6.48558742916689
$ws.Range("R9").Value = 58.37028686250201
$ws.Range("S9").Value = 0.1392157191583261
$ws.Range("T9").Value = 0.1392157191583261

$ws.Range("G10").Value = 0.3180153333333334
$ws.Range("H10").Value = 0.9540460000000001
$ws.Range("I10").Value = 0.2238320355298628
$ws.Range("J10").Value = 0.2238320355298628
$ws.Range("M10").Value = 9.887829999999999
$ws.Range("N10").Value = 29.66349
$ws.Range("O10").Value = 0.3015544655826039
$ws.Range("P10").Value = 0.301554465582604
$ws.Range("Q10").Value = 3.144481553393334
$ws.Range("R10").Value = 28.30033398054
$ws.Range("S10").Value = 0.06749754985447419
$ws.Range("T10").Value = 0.06749754985447419
